# Reorder the language/GDP-share table to descending order by value,
# and drop the "Swedish" and "Uzbek" rows that no longer appear in the
# finished dataset (table shrinks from A1:B23 to A1:B21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header: Language / 2010) stays as-is.

$ws.Range("A2").Value = "English"
$ws.Range("B2").Value = 22.66037303121328

$ws.Range("A3").Value = "Chinese"
$ws.Range("B3").Value = 15.05339433349864

$ws.Range("A4").Value = "Spanish"
$ws.Range("B4").Value = 7.119644881420998

$ws.Range("A5").Value = "Arabic"
$ws.Range("B5").Value = 5.278183174792866

$ws.Range("A6").Value = "Japanese"
$ws.Range("B6").Value = 5.029137285388825

$ws.Range("A7").Value = "German"
$ws.Range("B7").Value = 4.449831153007453

$ws.Range("A8").Value = "Portuguese"
$ws.Range("B8").Value = 3.615749093656355

$ws.Range("A9").Value = "Russian"
$ws.Range("B9").Value = 3.532555958236922

$ws.Range("A10").Value = "Malay-Indonesian"
$ws.Range("B10").Value = 2.932681142809074

$ws.Range("A11").Value = "French"
$ws.Range("B11").Value = 2.792939657568918

$ws.Range("A12").Value = "Italian"
$ws.Range("B12").Value = 2.344794373453612

$ws.Range("A13").Value = "Korean"
$ws.Range("B13").Value = 1.703257548539746

$ws.Range("A14").Value = "Persian"
$ws.Range("B14").Value = 1.585362524757927

$ws.Range("A15").Value = "Turkish"
$ws.Range("B15").Value = 1.420904013363344

$ws.Range("A16").Value = "Dutch"
$ws.Range("B16").Value = 1.323756251034521

$ws.Range("A17").Value = "Thai"
$ws.Range("B17").Value = 0.983578805990833

$ws.Range("A18").Value = "Polish"
$ws.Range("B18").Value = 0.8921610812276971

$ws.Range("A19").Value = "Urdu"
$ws.Range("B19").Value = 0.8641353632307527

$ws.Range("A20").Value = "Vietnamese"
$ws.Range("B20").Value = 0.5226989958802405

$ws.Range("A21").Value = "Bengali"
$ws.Range("B21").Value = 0.5145697962443675

# Remove the now-unused trailing rows (previously Swedish / Uzbek) so the
# sheet's used range / dimension shrinks to A1:B21.
$ws.Range("A22:B23").EntireRow.Delete()
